$d = $word.ActiveDocument

# The final paragraph in the document currently holds the (empty) "_GoBack"
# bookmark. Typing the new sentences there, separated by paragraph marks,
# grows 21 new paragraphs in front of it while the bookmark naturally stays
# attached to the trailing paragraph - exactly mirroring how this was done
# by hand in Word.
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertBefore("I look forward to hearing from you about these new obstacles`rI would be grateful if you take over the recap about the last sprint`rI would be grateful if you tell me how did you cope with that obstacle at the last sprint`rI would be grateful if you tell me about the new vacancy as soon as possible`rI would be grateful if you tell me about  the new due date of the current sprint`rI look forward to hearing from you about news from the morning meeting`rI would like to enquire about the new benchmark of the requirement’s integrity`rI would like to enquire about these new obstacles and how are we going to cope with them`rI would like to enquire about that unambiguous hint from the last meeting`rI would like to enquire about the new due date of the current sprint`rI would like to enquire about the new obstacles `rI would like to enquire about how are you going to cope with these new obstacles`rI would like to enquire the new due date`rI would like to enquire some information about the new benchmark of the requirement’s integrity`rI would like to enquire about that new vacancy`rI would like to enquire about the goals of this test cycle`rI would like to enquire about one obstacle in this test cycle`rI look forward to hearing from you about the new due date of the current sprint`rI look forward to hearing from you how do we suppose to cope with this obstacle`rI look forward to hearing from you about your decision regarding the new position`rI look forward to hearing from you about new goals in the next sprint`rI look forward to hearing from you about the new vacancy at the QA department")
